$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.450.66"
$ws.Range("E2").Value = "  -0.71%  "

$ws.Range("D3").Value = "3.216.17"
$ws.Range("E3").Value = "  -1.18%  "

$ws.Range("D5").Value = "'578.59"
$ws.Range("E5").Value = "  -0.88%  "

$ws.Range("D6").Value = "'182.74"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +0.77%  "

$ws.Range("D9").Value = "3.211.47"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("E10").Value = "  -3.30%  "

$ws.Range("E11").Value = "  -1.55%  "

$ws.Range("D12").Value = "'0.413"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").Value = "3.776.66"
$ws.Range("E13").Value = "  -1.18%  "

$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").Value = "'27.57"
$ws.Range("E15").Value = "  -3.26%  "

$ws.Range("D16").Value = "67.458.69"
$ws.Range("E16").Value = "  -0.70%  "

$ws.Range("E17").Value = "  -2.13%  "

$ws.Range("D18").Value = "3.205.58"
$ws.Range("E18").Value = "  -1.59%  "

$ws.Range("E19").Value = "  -2.06%  "

$ws.Range("D20").Value = "'13.40"
$ws.Range("E20").Value = "  -1.37%  "

$ws.Range("D21").Value = "'392.96"
$ws.Range("E21").Value = "  +2.93%  "

$ws.Range("E22").Value = "  -2.08%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'70.61"
$ws.Range("E24").Value = "  -1.04%  "

$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("E26").Value = "  -2.33%  "

$ws.Range("E27").Value = "  +1.35%  "

$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = "  -2.99%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("E30").Value = "  -1.95%  "

$ws.Range("E31").Value = "  -2.39%  "

$ws.Range("E32").Value = "  -1.56%  "

$ws.Range("D33").Value = "'6.95"
$ws.Range("E33").Value = "  -4.72%  "

$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("E35").Value = "  -1.56%  "

$ws.Range("D36").Value = "'161.40"
$ws.Range("E36").Value = "  -0.64%  "

$ws.Range("E37").Value = "  -5.73%  "

$ws.Range("E38").Value = "  +1.94%  "

$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'26.17"
$ws.Range("E39").Value = "  -1.68%  "

$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'0.801"
$ws.Range("E40").Value = "  -4.04%  "

$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").Value = "'6.45"
$ws.Range("E42").Value = "  -4.26%  "

$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "'0.0681"
$ws.Range("E43").Value = "  -0.91%  "

$ws.Range("D44").Value = "'40.62"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.44"
$ws.Range("E45").Value = "  -6.10%  "

$ws.Range("D46").Value = "2.602.85"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").Value = "'24.66"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "'333.97"
$ws.Range("E48").Value = "  -3.21%  "

$ws.Range("E49").Value = "  -2.98%  "

$ws.Range("E50").Value = "  +0.89%  "

$ws.Range("E51").Value = "  -1.86%  "
